$wb = $excel.ActiveWorkbook

# Update the PINI sheet values: A2 -> "a", B2 -> "1" (kept as text)
$ws = $wb.Worksheets.Item("PINI")
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "a"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"

# Delete the "Algo" sheet
$excel.DisplayAlerts = $false
$wsAlgo = $wb.Worksheets.Item("Algo")
$wsAlgo.Delete()
$excel.DisplayAlerts = $true
